$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.898.00'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -0.29%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.897.09'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -0.01%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.7919'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -4.33%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '243.77'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +0.82%  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.08%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3159'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -3.43%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '25.34'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -4.25%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07156'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +1.81%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08101'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +0.19%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.583'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +6.65%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.7671'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +1.07%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.869.25'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -1.58%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '92.56'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +0.50%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.162'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +5.42%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '29.917.75'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '13.93'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -1.09%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '244.33'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +0.19%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000007776'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +0.48%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '8.258'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +18.84%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.173.71'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +0.80%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.001'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -0.02%  '
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +0.00%  '
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -3.63%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.473'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +2.46%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '163.99'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -0.98%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.72'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -0.83%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.061'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -1.23%  '
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +2.90%  '
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +2.33%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.486'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +5.02%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05599'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -5.34%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.089'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +0.68%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.278'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +1.23%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7415'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +1.55%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.9996'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -0.08%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.634'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -3.33%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01931'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +0.88%  '
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +0.20%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.160.65'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +17.18%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '74.24'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +2.76%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.4421'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -0.25%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.948'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +1.69%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.8528'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +0.34%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '104.69'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +2.72%  '
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -0.10%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.881'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -0.56%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.963'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +1.93%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.454'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -1.07%  '
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +10.16%  '
